# Finetune homework: update the closing "thank you" slide's author credit.
#
# Slide 7 ("KÖSZÖNÖM A FIGYELMET!") has a text box with three paragraphs:
#   1. "KÖSZÖNÖM A FIGYELMET!"
#   2. "Kis Bertalan"
#   3. "Bertalan_Kis@epam.com"
#
# The edit collapses paragraphs 2 and 3 into a single paragraph reading
# "Farkas László" (an explicit white RGB color instead of the lt1 theme
# color), and removes the e-mail paragraph entirely.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item(1)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# Paragraph 2 is "Kis Bertalan". Deleting its own text plus its trailing
# paragraph mark merges it with paragraph 3 ("Bertalan_Kis@epam.com"),
# which inherits paragraph 3's paragraph properties (no spcAft) while
# paragraph 3's run formatting (italic Calibri) is still sitting there
# ready to be overwritten.
$para2 = $tr.Paragraphs(2, 1)
$mergeRange = $tr.Characters($para2.Start, $para2.Length + 1)
$mergeRange.Delete()

# The merged paragraph (still index 2) now holds the leftover text from
# paragraph 3; replace it with the new name and force an explicit white
# fill color (rather than the lt1 scheme color it had before).
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "Farkas László"
$para2.Font.Color.RGB = 16777215
